$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header banners: new report date 25.12.2024 (previously 18.12.2024)
$ws.Range("A1").Value = "Mangrove Communication   25.12.2024"
$ws.Range("A10").Value = "DAILY STOCK                         (25/12/2024) "

# --- RSO sale table (rows 3-6) ---
$ws.Range("C3").Value = 31645
$ws.Range("C4").Value = 22199
$ws.Range("F4").Value = 25
$ws.Range("C5").Value = 24831
$ws.Range("F5").Value = 5
$ws.Range("C6").Value = 36442

# --- Daily stock table ---
$ws.Range("C13").Value = 190819

$ws.Range("C14").Value = 375420
$ws.Range("D14").Value = 115117
$ws.Range("E14").ClearContents()

$ws.Range("C18").Value = 500

$ws.Range("C20").Value = 2000

$ws.Range("E21").ClearContents()

$ws.Range("C24").Value = 25

$ws.Range("C25").Value = 39
$ws.Range("D25").ClearContents()

$ws.Range("D26").Value = 30

$ws.Range("C27").Value = 77
$ws.Range("D27").Value = 2

# Update current selection to match the saved workbook state
$ws.Range("J31").Select()
